$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 5853.5713
$ws.Range("J69").Value = 5853.5713
$ws.Range("L69").Value = 17560.7139
$ws.Range("N69").Value = -19308.7139
$ws.Range("H72").Value = 5853.5713
$ws.Range("J72").Value = 5853.5713
$ws.Range("L72").Value = 52682.14169999999
$ws.Range("N72").Value = -61418.14169999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3694.762
$ws.Range("I2").Value = 1229.2858
$ws.Range("J2").Value = 8625.714
$ws.Range("K2").Value = 1229.2858
$ws.Range("L2").Value = 8625.714
$ws.Range("M2").Value = -1116.2858
$ws.Range("N2").Value = -8851.714
$ws.Range("H4").Value = 349.25
$ws.Range("I4").Value = 342
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 342
$ws.Range("L4").Value = 400
$ws.Range("M4").Value = -226
$ws.Range("N4").Value = -632
$ws.Range("H5").Value = 342.75
$ws.Range("I5").Value = 356.33334
$ws.Range("J5").Value = 302
$ws.Range("K5").Value = 356.33334
$ws.Range("L5").Value = 302
$ws.Range("M5").Value = -244.33334
$ws.Range("N5").Value = -526
$ws.Range("H9").Value = 17777
$ws.Range("J9").Value = 17777
$ws.Range("L9").Value = 17777
$ws.Range("N9").Value = -18117
$ws.Range("H20").Value = 17777
$ws.Range("J20").Value = 17777
$ws.Range("L20").Value = 17777
$ws.Range("N20").Value = -18317
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H37").Value = 10120.667
$ws.Range("J37").Value = 9508.667
$ws.Range("L37").Value = 9508.667
$ws.Range("N37").Value = -10054.667
$ws.Range("H44").Value = 37719.5
$ws.Range("J44").Value = 37719.5
$ws.Range("L44").Value = 37719.5
$ws.Range("N44").Value = -38695.5
$ws.Range("H55").Value = 24042.5
$ws.Range("J55").Value = 24042.5
$ws.Range("L55").Value = 24042.5
$ws.Range("N55").Value = -24672.5
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H116").Value = 3694.762
$ws.Range("I116").Value = 1229.2858
$ws.Range("J116").Value = 8625.714
$ws.Range("K116").Value = 1229.2858
$ws.Range("L116").Value = 8625.714
$ws.Range("M116").Value = 1064.7142
$ws.Range("N116").Value = -13213.714

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3694.762
$ws.Range("I3").Value = 1229.2858
$ws.Range("J3").Value = 8625.714
$ws.Range("K3").Value = 1229.2858
$ws.Range("L3").Value = 8625.714
$ws.Range("M3").Value = -1115.2858
$ws.Range("N3").Value = -8853.714
$ws.Range("H4").Value = 342.75
$ws.Range("I4").Value = 356.33334
$ws.Range("J4").Value = 302
$ws.Range("K4").Value = 356.33334
$ws.Range("L4").Value = 302
$ws.Range("M4").Value = -241.33334
$ws.Range("N4").Value = -532
$ws.Range("H15").Value = 40000
$ws.Range("J15").Value = 40000
$ws.Range("L15").Value = 40000
$ws.Range("N15").Value = -40454
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()
$ws.Range("H22").Value = 599.5
$ws.Range("I22").Value = 499
$ws.Range("J22").Value = 700
$ws.Range("K22").Value = 499
$ws.Range("L22").Value = 700
$ws.Range("M22").Value = -326
$ws.Range("N22").Value = -1046
$ws.Range("H35").Value = 34887
$ws.Range("J35").Value = 34887
$ws.Range("L35").Value = 34887
$ws.Range("N35").Value = -35507
$ws.Range("H82").Value = 7128.5
$ws.Range("I82").Value = 7128.5
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 7128.5
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -6745.5
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 7128.5
$ws.Range("I85").Value = 7128.5
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 7128.5
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -5802.5
$ws.Range("N85").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 35295
$ws.Range("J68").Value = 35295
$ws.Range("L68").Value = 35295
$ws.Range("N68").Value = -36793
$ws.Range("H71").Value = 35295
$ws.Range("J71").Value = 35295
$ws.Range("L71").Value = 105885
$ws.Range("N71").Value = -113373
$ws.Range("H74").Value = 35314
$ws.Range("J74").Value = 35314
$ws.Range("L74").Value = 35314
$ws.Range("N74").Value = -37062
$ws.Range("H77").Value = 35314
$ws.Range("J77").Value = 35314
$ws.Range("L77").Value = 105942
$ws.Range("N77").Value = -114678

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3406.1177
$ws.Range("I80").Value = 3077.2307
$ws.Range("K80").Value = 3077.2307
$ws.Range("M80").Value = -2079.2307
$ws.Range("H83").Value = 3406.1177
$ws.Range("I83").Value = 3077.2307
$ws.Range("K83").Value = 15386.1535
$ws.Range("M83").Value = -10394.1535
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1867.3334
$ws.Range("I68").Value = 1867.3334
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1867.3334
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -1118.3334
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 1867.3334
$ws.Range("I71").Value = 1867.3334
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 9336.667
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -5592.666999999999
$ws.Range("N71").ClearContents()
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7001.1665
$ws.Range("I81").Value = 13252.25
$ws.Range("J81").Value = 2000.3
$ws.Range("K81").Value = 26504.5
$ws.Range("L81").Value = 4000.6
$ws.Range("M81").Value = -25443.5
$ws.Range("N81").Value = -6122.6
$ws.Range("H84").Value = 7001.1665
$ws.Range("I84").Value = 13252.25
$ws.Range("J84").Value = 2000.3
$ws.Range("K84").Value = 132522.5
$ws.Range("L84").Value = 20003
$ws.Range("M84").Value = -127218.5
$ws.Range("N84").Value = -30611
